$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value2 = 111909174
$ws.Range("B15").Value2 = 77267
$ws.Range("D15").Value2 = 'NT'
$ws.Range("E15").Value2 = 6446
$ws.Range("F15").Value2 = 'Kolflarnlav'
$ws.Range("G15").Value2 = 'Carbonicola anthracophila'
$ws.Range("H15").Value2 = '(Nyl.) Bendiksby & Timdal'
$ws.Range("Q15").Value2 = 467989.0228066717
$ws.Range("R15").Value2 = 6875352.744105402
# Row 17
$ws.Range("A17").Value2 = 111908768
$ws.Range("B17").Value2 = 96348
$ws.Range("D17").Value2 = 'VU'
$ws.Range("E17").Value2 = 220787
$ws.Range("F17").Value2 = 'Knärot'
$ws.Range("G17").Value2 = 'Goodyera repens'
$ws.Range("H17").Value2 = '(L.) R. Br.'
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value2 = '1'
$ws.Range("P17").Value2 = 'Fläcksberget, Hjd'
$ws.Range("Q17").Value2 = 467911.8445363804
$ws.Range("R17").Value2 = 6875299.456096188
# Row 18
$ws.Range("A18").Value2 = 111908364
$ws.Range("B18").Value2 = 90660
$ws.Range("E18").Value2 = 4362
$ws.Range("F18").Value2 = 'Blå taggsvamp'
$ws.Range("G18").Value2 = 'Hydnellum caeruleum'
$ws.Range("H18").Value2 = '(Hornem.) P.Karst.'
$ws.Range("P18").Value2 = 'Gröbäcken, Hjd'
$ws.Range("Q18").Value2 = 467724.2196293612
$ws.Range("R18").Value2 = 6874811.291555981
# Row 19
$ws.Range("A19").Value2 = 111909766
$ws.Range("B19").Value2 = 89183
$ws.Range("D19").Value2 = 'LC'
$ws.Range("E19").Value2 = 3215
$ws.Range("F19").Value2 = 'Rödgul trumpetsvamp'
$ws.Range("G19").Value2 = 'Craterellus lutescens'
$ws.Range("H19").Value2 = '(Fr.) Fr.'
$ws.Range("I19").Value2 = ""
$ws.Range("Q19").Value2 = 467756.8135427741
$ws.Range("R19").Value2 = 6875469.545251801
# Row 21
$ws.Range("A21").Value2 = 112014423
$ws.Range("B21").Value2 = 90658
$ws.Range("E21").Value2 = 4361
$ws.Range("F21").Value2 = 'Orange taggsvamp'
$ws.Range("G21").Value2 = 'Hydnellum aurantiacum'
$ws.Range("H21").Value2 = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q21").Value2 = 467430.0274016621
$ws.Range("R21").Value2 = 6875237.811246304
# Row 22
$ws.Range("A22").Value2 = 112014347
$ws.Range("B22").Value2 = 90678
$ws.Range("E22").Value2 = 4366
$ws.Range("F22").Value2 = 'Skarp dropptaggsvamp'
$ws.Range("G22").Value2 = 'Hydnellum peckii'
$ws.Range("H22").Value2 = 'Banker'
$ws.Range("Q22").Value2 = 467430.0274016621
$ws.Range("R22").Value2 = 6875237.811246304
# Row 23
$ws.Range("A23").Value2 = 112014300
$ws.Range("B23").Value2 = 90689
$ws.Range("E23").Value2 = 5966
$ws.Range("F23").Value2 = 'Motaggsvamp'
$ws.Range("G23").Value2 = 'Sarcodon squamosus'
$ws.Range("H23").Value2 = '(Schaeff.) Quél.'
$ws.Range("Q23").Value2 = 467415.4484496959
$ws.Range("R23").Value2 = 6875287.271149865
# Row 24
$ws.Range("A24").Value2 = 112015011
$ws.Range("B24").Value2 = 90658
$ws.Range("E24").Value2 = 4361
$ws.Range("F24").Value2 = 'Orange taggsvamp'
$ws.Range("G24").Value2 = 'Hydnellum aurantiacum'
$ws.Range("H24").Value2 = '(Batsch:Fr.) P.Karst.'
# Row 25
$ws.Range("A25").Value2 = 112014177
$ws.Range("Q25").Value2 = 467389.9660160011
$ws.Range("R25").Value2 = 6875327.91063729
# Row 26
$ws.Range("A26").Value2 = 112014142
$ws.Range("B26").Value2 = 90666
$ws.Range("E26").Value2 = 4364
$ws.Range("F26").Value2 = 'Dropptaggsvamp'
$ws.Range("G26").Value2 = 'Hydnellum ferrugineum'
$ws.Range("H26").Value2 = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q26").Value2 = 467442.7363991642
$ws.Range("R26").Value2 = 6875336.798642672
# Row 27
$ws.Range("A27").Value2 = 112014923
$ws.Range("B27").Value2 = 90689
$ws.Range("E27").Value2 = 5966
$ws.Range("F27").Value2 = 'Motaggsvamp'
$ws.Range("G27").Value2 = 'Sarcodon squamosus'
$ws.Range("H27").Value2 = '(Schaeff.) Quél.'
$ws.Range("Q27").Value2 = 467413.0579403224
$ws.Range("R27").Value2 = 6875234.216212902
# Row 28
$ws.Range("A28").Value2 = 112014208
$ws.Range("Q28").Value2 = 467418.043506761
$ws.Range("R28").Value2 = 6875312.610613029
# Row 29
$ws.Range("A29").Value2 = 112014229
$ws.Range("B29").Value2 = 90682
$ws.Range("E29").Value2 = 2059
$ws.Range("F29").Value2 = 'Skrovlig taggsvamp'
$ws.Range("G29").Value2 = 'Hydnellum scabrosum'
$ws.Range("H29").Value2 = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q29").Value2 = 467427.230114766
$ws.Range("R29").Value2 = 6875289.506732536
